$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.400352648228136
$ws.Range("C2").Value = 0.2522079696473725
$ws.Range("D2").Value = 0.09410523387798975
$ws.Range("F2").Value = 2.175727196881596
$ws.Range("G2").Value = 0.002511508971311874
$ws.Range("I2").Value = 1.649109537196878
$ws.Range("L2").Value = 0.247266623261801
$ws.Range("M2").Value = 0.3093111431223861
$ws.Range("B3").Value = 1.313779732063551
$ws.Range("C3").Value = 0.2192766030679252
$ws.Range("D3").Value = 0.09418093707741093
$ws.Range("F3").Value = 2.120018071462468
$ws.Range("G3").Value = 0.00251709033857021
$ws.Range("I3").Value = 1.62632282151101
$ws.Range("L3").Value = 0.2451822694918135
$ws.Range("M3").Value = 0.2965791635638411
$ws.Range("B4").Value = 1.261534399895652
$ws.Range("C4").Value = 0.1990517184434566
$ws.Range("D4").Value = 0.09425349622320667
$ws.Range("F4").Value = 2.087038882756602
$ws.Range("G4").Value = 0.002520696324789328
$ws.Range("I4").Value = 1.613122399793212
$ws.Range("L4").Value = 0.244015803945409
$ws.Range("M4").Value = 0.2889608564446817
$ws.Range("B5").Value = 1.240472312984366
$ws.Range("C5").Value = 0.1908079871277266
$ws.Range("D5").Value = 0.09428956836805469
$ws.Range("F5").Value = 2.073905552361438
$ws.Range("G5").Value = 0.002522210963733091
$ws.Range("I5").Value = 1.607940564846885
$ws.Range("L5").Value = 0.2435689273792647
$ws.Range("M5").Value = 0.2859063091627547
$ws.Range("B6").Value = 1.236988738825119
$ws.Range("C6").Value = 0.1894389853224254
$ws.Range("D6").Value = 0.09429594936168328
$ws.Range("F6").Value = 2.071743180806152
$ws.Range("G6").Value = 0.00252246520149353
$ws.Range("I6").Value = 1.607092009778214
$ws.Range("L6").Value = 0.2434964423440675
$ws.Range("M6").Value = 0.2854021191604019
$ws.Range("B7").Value = 1.261249425732501
$ws.Range("C7").Value = 0.1989405493055472
$ws.Range("D7").Value = 0.09425395644245071
$ws.Range("F7").Value = 2.086860526506285
$ws.Range("G7").Value = 0.002520716568844522
$ws.Range("I7").Value = 1.613051718163078
$ws.Range("L7").Value = 0.2440096619874197
$ws.Range("M7").Value = 0.2889194595492839
$ws.Range("B8").Value = 1.370312800604268
$ws.Range("C8").Value = 0.2408538186309102
$ws.Range("D8").Value = 0.09412589170592867
$ws.Range("F8").Value = 2.156262649042844
$ws.Range("G8").Value = 0.002513396370092913
$ws.Range("I8").Value = 1.641087719090066
$ws.Range("L8").Value = 0.2465243945721411
$ws.Range("M8").Value = 0.3048797436943644
$ws.Range("B9").Value = 1.591459354116864
$ws.Range("C9").Value = 0.3230423801514348
$ws.Range("D9").Value = 0.09408410990698712
$ws.Range("F9").Value = 2.302210485848633
$ws.Range("G9").Value = 0.002500454565711556
$ws.Range("I9").Value = 1.702406926146352
$ws.Range("L9").Value = 0.2523570612463075
$ws.Range("M9").Value = 0.3377655279360567
$ws.Range("B10").Value = 1.758456058550109
$ws.Range("C10").Value = 0.3834796468014474
$ws.Range("D10").Value = 0.09418445972587364
$ws.Range("F10").Value = 2.41562175383271
$ws.Range("G10").Value = 0.002491797524859339
$ws.Range("I10").Value = 1.751422018503149
$ws.Range("L10").Value = 0.2571954013645978
$ws.Range("M10").Value = 0.3629081114560577
$ws.Range("B11").Value = 1.835429254279802
$ws.Range("C11").Value = 0.4109981571328944
$ws.Range("D11").Value = 0.09425931019356426
$ws.Range("F11").Value = 2.468597651429036
$ws.Range("G11").Value = 0.002488041895268189
$ws.Range("I11").Value = 1.774602669743146
$ws.Range("L11").Value = 0.2595174422384474
$ws.Range("M11").Value = 0.3745624383629149
$ws.Range("B12").Value = 1.864722921158773
$ws.Range("C12").Value = 0.4214232432694303
$ws.Range("D12").Value = 0.09429191352460009
$ws.Range("F12").Value = 2.48886031353274
$ws.Range("G12").Value = 0.002486645813290167
$ws.Range("I12").Value = 1.783509271650487
$ws.Range("L12").Value = 0.2604142026841885
$ws.Range("M12").Value = 0.3790070164975248
$ws.Range("B13").Value = 1.858407514518717
$ws.Range("C13").Value = 0.4191778050329162
$ws.Range("D13").Value = 0.09428470142916012
$ws.Range("F13").Value = 2.484487365842796
$ws.Range("G13").Value = 0.002486945326433872
$ws.Range("I13").Value = 1.781585331487818
$ws.Range("L13").Value = 0.2602202922534929
$ws.Range("M13").Value = 0.3780484007665166
$ws.Range("B14").Value = 1.837836338239697
$ws.Range("C14").Value = 0.4118557420333673
$ws.Range("D14").Value = 0.09426190678551905
$ws.Range("F14").Value = 2.470260611856787
$ws.Range("G14").Value = 0.002487926516656191
$ws.Range("I14").Value = 1.775332835108543
$ws.Range("L14").Value = 0.2595908691856152
$ws.Range("M14").Value = 0.374927467511263
$ws.Range("B15").Value = 1.825254888862389
$ws.Range("C15").Value = 0.4073713666814456
$ws.Range("D15").Value = 0.09424850085460434
$ws.Range("F15").Value = 2.461572681612097
$ws.Range("G15").Value = 0.002488530917774904
$ws.Range("I15").Value = 1.771519793425227
$ws.Range("L15").Value = 0.2592076038042848
$ws.Range("M15").Value = 0.3730198915986023
$ws.Range("B16").Value = 1.753445997150379
$ws.Range("C16").Value = 0.381681823299175
$ws.Range("D16").Value = 0.09418016079111169
$ws.Range("F16").Value = 2.412187748793542
$ws.Range("G16").Value = 0.002492046623909466
$ws.Range("I16").Value = 1.749925018111114
$ws.Range("L16").Value = 0.2570460904050123
$ws.Range("M16").Value = 0.3621508503406758
$ws.Range("B17").Value = 1.70965180510791
$ws.Range("C17").Value = 0.3659291538778575
$ws.Range("D17").Value = 0.09414575799872438
$ws.Range("F17").Value = 2.382248250240167
$ws.Range("G17").Value = 0.002494250031702801
$ws.Range("I17").Value = 1.736904695258247
$ws.Range("L17").Value = 0.2557511096630805
$ws.Range("M17").Value = 0.3555386922795591
$ws.Range("B18").Value = 1.68455712986713
$ws.Range("C18").Value = 0.356870968014789
$ws.Range("D18").Value = 0.09412871533110234
$ws.Range("F18").Value = 2.365157996587044
$ws.Range("G18").Value = 0.002495534560403062
$ws.Range("I18").Value = 1.729498809771016
$ws.Range("L18").Value = 0.2550176624637857
$ws.Range("M18").Value = 0.3517559540364488
$ws.Range("B19").Value = 1.676076719789762
$ws.Range("C19").Value = 0.3538043988961022
$ws.Range("D19").Value = 0.09412341459392337
$ws.Range("F19").Value = 2.359393801268112
$ws.Range("G19").Value = 0.002495972436260763
$ws.Range("I19").Value = 1.727005522626214
$ws.Range("L19").Value = 0.2547712849550265
$ws.Range("M19").Value = 0.350478682971314
$ws.Range("B20").Value = 1.714303975328448
$ws.Range("C20").Value = 0.3676058053295606
$ws.Range("D20").Value = 0.09414913571469441
$ws.Range("F20").Value = 2.385421870739179
$ws.Range("G20").Value = 0.00249401369737936
$ws.Range("I20").Value = 1.73828212325607
$ws.Range("L20").Value = 0.2558877832690456
$ws.Range("M20").Value = 0.3562404553424443
$ws.Range("B21").Value = 1.843874632017958
$ws.Range("C21").Value = 0.4140062830360307
$ws.Range("D21").Value = 0.09426848605414051
$ws.Range("F21").Value = 2.474433854492219
$ws.Range("G21").Value = 0.002487637610428944
$ws.Range("I21").Value = 1.777165842441121
$ws.Range("L21").Value = 0.2597752719255766
$ws.Range("M21").Value = 0.3758433094906337
$ws.Range("B22").Value = 1.929405638982018
$ws.Range("C22").Value = 0.4443578364426344
$ws.Range("D22").Value = 0.09437134176632611
$ws.Range("F22").Value = 2.533786099696869
$ws.Range("G22").Value = 0.002483622492742615
$ws.Range("I22").Value = 1.803328866411277
$ws.Range("L22").Value = 0.262417728153352
$ws.Range("M22").Value = 0.3888376442214181
$ws.Range("B23").Value = 1.883678099376482
$ws.Range("C23").Value = 0.4281559840623004
$ws.Range("D23").Value = 0.094314151705575
$ws.Range("F23").Value = 2.502000009554763
$ws.Range("G23").Value = 0.002485751575664559
$ws.Range("I23").Value = 1.789296001610637
$ws.Range("L23").Value = 0.2609980729910433
$ws.Range("M23").Value = 0.3818855493577331
$ws.Range("B24").Value = 1.712200470130199
$ws.Range("C24").Value = 0.3668477967688659
$ws.Range("D24").Value = 0.09414760013189749
$ws.Range("F24").Value = 2.383986695985556
$ws.Range("G24").Value = 0.00249412048883511
$ws.Range("I24").Value = 1.737659140016518
$ws.Range("L24").Value = 0.2558259587032126
$ws.Range("M24").Value = 0.3559231300908721
$ws.Range("B25").Value = 1.530845922881099
$ws.Range("C25").Value = 0.300802875935176
$ws.Range("D25").Value = 0.09407264808275784
$ws.Range("F25").Value = 2.261654829178156
$ws.Range("G25").Value = 0.002503805437852961
$ws.Range("I25").Value = 1.685129508552123
$ws.Range("L25").Value = 0.2506823115701451
$ws.Range("M25").Value = 0.2859063091627547
